$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final dataset for rows 2..31 (A..H) after the May 9th edit: 11 new sensor
# readings were prepended (rows 2-12) and the previously-existing readings
# were kept (now rows 13-31, timestamps renumbered to continue the sequence),
# dropping the final old row that no longer fits the 30-row window.
$data = @(
    @(0,"falling",-4.1851313829422,4.315697193145752,-1.521726083755494,0.0226020142436027,0.0108428578823804,0.0673478916287422),
    @(100,"falling",-4.044641017913818,4.278119564056396,-1.635269999504089,-0.0668897479772567,0.0363464802503585,0.0458148941397666),
    @(200,"falling",-4.037833023071289,4.246780538558959,-1.5330271422863,-0.0578794814646244,-0.0059559359215199,-0.026419922709465),
    @(300,"falling",-4.188077163696289,4.206779670715332,-1.555011582374573,-0.0319177098572254,0.0326812900602817,0.0067195175215601),
    @(400,"falling",-4.266698265075684,4.172728943824768,-1.657637816667557,-0.0120645882561802,-0.0032070425804704,0.0010690141934901),
    @(500,"falling",-4.307723093032837,4.115167903900146,-1.525768172740936,0.00534507073462,-0.0042760567739605,0.0001527163112768),
    @(600,"falling",-4.275023818016052,4.170393228530884,-1.397946566343308,0.0019853119738399,-0.0022907445672899,0.0009162978967651),
    @(700,"falling",-4.101493644714355,4.221428394317627,-1.457295727729798,0.009010262787342,-0.0284052342176437,-0.0154243474826216),
    @(800,"falling",-4.052593421936035,4.150938177108764,-1.483953177928925,0.027030786499381,0.0221438650041818,-0.0340557359158992),
    @(900,"falling",-4.146867895126342,4.078793859481812,-1.569497275352478,0.1010981947183609,-0.0036651915870606,-0.0352774672210216),
    @(1000,"falling",-4.033952045440674,4.105497026443481,-1.389435482025145,-0.0256563406437635,0.0059559359215199,-0.0114537235349416),
    @(1100,"falling",-3.75560998916626,3.84902286529541,-0.6511768102645874,-0.0073303831741213,-0.0361937656998634,-0.0346666015684604),
    @(1200,"falling",-4.870891666412356,3.396544885635375,-0.1143038809299458,-0.00335975876078,-0.1609629988670349,0.08643743395805351),
    @(1300,"falling",-4.700507307052611,4.103157997131351,-0.4222948729991927,0.08429940044879911,0.1406517177820205,-0.2235766798257827),
    @(1400,"falling",-4.943684387207032,4.186895227432251,0.9878155738115346,0.4598288238048553,0.2849686443805694,-0.0438295826315879),
    @(1500,"falling",-5.002620124816895,2.714760780334467,1.064968466758727,-0.0018325957935303,0.4227187633514404,-0.1143845170736312),
    @(1600,"falling",-5.456374883651734,2.089365065097808,1.750832557678223,0.0514653958380222,0.7043276429176331,-0.5734497308731079),
    @(1700,"falling",-5.036787080764768,1.541897773742674,1.896076297760009,0.2469422817230224,-1.195616006851196,-1.284649610519409),
    @(1800,"falling",-5.012954092025758,0.4318521261215144,0.7947881996631543,0.1327104717493057,-1.82404363155365,-1.386511445045471),
    @(1900,"falling",-6.129836273193367,1.83154973983766,-1.60779109001161,-1.465160250663757,4.381736278533936,-4.054465293884277),
    @(2000,"falling",-4.299563241004924,7.852099323272741,-5.872165632247949,2.828611612319946,0.6840163469314575,1.775785326957703),
    @(2100,"falling",-3.646280765533453,12.82203674316406,-5.478024601936328,-0.8000807762145996,0.6334672570228577,-1.410640597343445),
    @(2200,"falling",-2.318737292289721,6.985545444488493,-2.094196039438245,0.380110889673233,0.3691153228282928,0.4048509299755096),
    @(2300,"falling",1.412541198730471,1.629453468322767,-1.321827816963195,-0.2186897546052932,0.1003346145153045,0.2921462953090668),
    @(2400,"falling",1.654833269119258,6.014033746719358,-2.261167562007916,-0.1091921627521514,-0.4050036668777466,0.2264782935380935),
    @(2500,"falling",-0.3583765506744629,5.71919040679932,-3.682787048816667,-0.4273002445697784,0.8710938692092896,0.479834645986557),
    @(2600,"falling",-2.107488989830001,6.458877801895146,-2.10113409161568,-0.0183259565383195,0.1708895564079284,0.3941608071327209),
    @(2700,"falling",2.009711074829105,7.251273727416988,-2.904031562805179,-0.4156937897205353,-0.4030183553695678,-0.1386664062738418),
    @(2800,"falling",2.518718385696409,6.290945339202878,-3.275830793380734,-0.2107485085725784,-0.07910704612731929,0.06582073122262951),
    @(2900,"falling",1.816467857360826,5.365327644348143,-2.097930765151977,-0.2503020465373993,0.09987647086381909,0.2165517359972)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}
